$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.712.14"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").Value = "2.209.18"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'254.70"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").Value = "'75.06"
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.588"
$ws.Range("E9").Value = "  -4.47%  "
$ws.Range("D10").Value = "'41.13"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").Value = "'0.0911"
$ws.Range("E11").Value = "  -2.44%  "
$ws.Range("D12").Value = "'6.87"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").Value = "2.540.43"
$ws.Range("D15").Value = "'14.29"
$ws.Range("E15").Value = "  -2.41%  "
$ws.Range("D16").Value = "2.215.37"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").Value = "'0.781"
$ws.Range("E17").Value = "  -3.69%  "
$ws.Range("D18").Value = "42.639.62"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").Value = "'0.0000102"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").Value = "'70.97"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "'5.92"
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").Value = "'228.68"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").Value = "'9.33"
$ws.Range("E24").Value = "  -8.90%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("D28").Value = "'39.18"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("E29").Value = "  -2.95%  "
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").Value = "'173.01"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "'20.19"
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").Value = "'0.0848"
$ws.Range("E33").Value = "  +6.72%  "
$ws.Range("D34").Value = "'5.19"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("E37").Value = "  +5.67%  "
$ws.Range("D38").Value = "'4.29"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "'12.30"
$ws.Range("E39").Value = "  -4.45%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("D41").Value = "'2.76"
$ws.Range("E41").Value = "  +18.85%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.197"
$ws.Range("E42").Value = "  -3.22%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'59.81"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").Value = "'5.23"
$ws.Range("E44").Value = "  -6.27%  "
$ws.Range("D45").Value = "'101.41"
$ws.Range("E45").Value = "  -5.49%  "
$ws.Range("D46").Value = "'8.35"
$ws.Range("D47").Value = "'0.0975"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("D48").Value = "'0.459"
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.66"
$ws.Range("E51").Value = "  -0.93%  "
